# Update "Epoch Accuracy" worksheet B2:B118 values and A102:A118 object repr text
# to match the new epoch-accuracy run, per commit "M08 Froze TE+D123 Data Aug".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.421875, 0.265625, 0.21875, 0.203125, 0.234375, 0.25, 0.203125, 0.15625, 0.1875, 0.140625,
    0.171875, 0.140625, 0.234375, 0.125, 0.125, 0.203125, 0.21875, 0.171875, 0.171875, 0.140625,
    0.125, 0.109375, 0.15625, 0.140625, 0.234375, 0.1875, 0.234375, 0.125, 0.109375, 0.15625,
    0.171875, 0.125, 0.15625, 0.109375, 0.078125, 0.1875, 0.15625, 0.15625, 0.203125, 0.15625,
    0.171875, 0.125, 0.09375, 0.125, 0.15625, 0.109375, 0.125, 0.140625, 0.140625, 0.125,
    0.125, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375,
    0.09375, 0.09375, 0.09375, 0.109375, 0.109375, 0.109375, 0.09375, 0.09375, 0.09375, 0.09375,
    0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.078125, 0.078125, 0.078125, 0.078125,
    0.078125, 0.078125, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375,
    0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375, 0.09375,
    0.09375, 0.21875, 0.0625, 0.125, 0.171875, 0.046875, 0.109375, 0.0625, 0.140625, 0.125,
    0.015625, 0.140625, 0.109375, 0.1875, 0.125, 0.125, 0.1147540983606557
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = [double]$values[$i]
}

$newAddr = "<__main__.DisplayOutputs object at 0x7ff7b862b9d0>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $newAddr
}
